# "final jing jing la"
# Populate the Receipt sheet with the new transaction rows (2-17).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row, date serial, description, amount column (C or D), amount
$entries = @(
    @{Row=2;  Date=45063; Desc="SCB Easy";         Col="D"; Amt=1720},
    @{Row=3;  Date=45062; Desc="SCB Easy";         Col="D"; Amt=1020},
    @{Row=4;  Date=45062; Desc="SCB Easy";         Col="D"; Amt=68},
    @{Row=5;  Date=45061; Desc="SCB Easy";         Col="D"; Amt=183},
    @{Row=6;  Date=45061; Desc="SCB Easy";         Col="D"; Amt=140},
    @{Row=7;  Date=45061; Desc="SCB Easy";         Col="D"; Amt=1000},
    @{Row=8;  Date=45059; Desc="SCB Easy";         Col="D"; Amt=73},
    @{Row=9;  Date=45058; Desc="SCB Easy";         Col="D"; Amt=45},
    @{Row=10; Date=45058; Desc="SCB Easy";         Col="D"; Amt=18},
    @{Row=11; Date=45058; Desc="SCB Easy";         Col="C"; Amt=340.38},
    @{Row=12; Date=45058; Desc="SCB Easy";         Col="C"; Amt=320.38},
    @{Row=13; Date=45058; Desc="SCB Easy";         Col="D"; Amt=1292},
    @{Row=14; Date=45055; Desc="SCB Easy";         Col="D"; Amt=141},
    @{Row=15; Date=57242; Desc="7-Eleven";         Col="D"; Amt=32},
    @{Row=16; Date=43629; Desc="The Sand Dollar";  Col="D"; Amt=82.98},
    @{Row=17; Date=43789; Desc="HARBOR LANE CAFE"; Col="D"; Amt=31.39}
)

foreach ($e in $entries) {
    $r = $e.Row
    $ws.Cells.Item($r, 1).Value = $e.Date
    $ws.Cells.Item($r, 2).Value = $e.Desc
    $ws.Range("$($e.Col)$r").Value = $e.Amt
}

# Column A holds dates -> give it a date number format. Set it on A2 first,
# then copy that cell's format onto the rest of the column so every cell
# shares a single style record (instead of minting one xf per cell).
$ws.Range("A2").NumberFormat = "mm-dd-yy"
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A3:A17").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

# Match the final selection left behind in the saved file
$ws.Range("G10").Select() | Out-Null
